# Update the yearly income-statement figures (rial.xlsx) with the
# latest published financial data, and clear the placeholder "-" cells
# in rows 15 and 23 back to numeric 0 now that real figures exist for
# the rest of the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Overview")

# Row 11: فروش (Sales)
$ws.Range("D11").Value = 2165020
$ws.Range("E11").Value = 2859959
$ws.Range("F11").Value = 4233243
$ws.Range("G11").Value = 6211749
$ws.Range("H11").Value = 14977487

# Row 12: بهای تمام شده کالای فروش رفته (Cost of goods sold)
$ws.Range("D12").Value = -1924597
$ws.Range("E12").Value = -2527347
$ws.Range("F12").Value = -3624067
$ws.Range("G12").Value = -5415675
$ws.Range("H12").Value = -12683135

# Row 13: سود (زیان) ناخالص (Gross profit)
$ws.Range("D13").Value = 240423
$ws.Range("E13").Value = 332612
$ws.Range("F13").Value = 609176
$ws.Range("G13").Value = 796074
$ws.Range("H13").Value = 2294352

# Row 14: هزینه های عمومی, اداری و تشکیلاتی (G&A expenses)
$ws.Range("D14").Value = -91425
$ws.Range("E14").Value = -111606
$ws.Range("F14").Value = -191409
$ws.Range("G14").Value = -265580
$ws.Range("H14").Value = -499564

# Row 15: هزینه کاهش ارزش دریافتنی‌‏ها (هزینه استثنایی) -> now 0 instead of "-"
$ws.Range("D15").Value = 0
$ws.Range("E15").Value = 0
$ws.Range("F15").Value = 0
$ws.Range("G15").Value = 0
$ws.Range("H15").Value = 0

# Row 16: خالص سایر درامدها (هزینه ها) ی عملیاتی
$ws.Range("D16").Value = -9444
$ws.Range("E16").Value = -12992
$ws.Range("F16").Value = -14806
$ws.Range("G16").Value = -22779
$ws.Range("H16").Value = -48155

# Row 17: سود (زیان) عملیاتی
$ws.Range("D17").Value = 139554
$ws.Range("E17").Value = 208014
$ws.Range("F17").Value = 402961
$ws.Range("G17").Value = 507715
$ws.Range("H17").Value = 1746633

# Row 18: هزینه های مالی
$ws.Range("D18").Value = -18739
$ws.Range("E18").Value = -20996
$ws.Range("F18").Value = -31043
$ws.Range("G18").Value = -69136
$ws.Range("H18").Value = -103359

# Row 19: خالص سایر درامدها و هزینه های غیرعملیاتی
$ws.Range("D19").Value = 18494
$ws.Range("E19").Value = 4402
$ws.Range("F19").Value = 7215
$ws.Range("G19").Value = -22393
$ws.Range("H19").Value = 18429

# Row 20: سود (زیان) خالص عملیات در حال تداوم قبل از مالیات
$ws.Range("D20").Value = 139309
$ws.Range("E20").Value = 191420
$ws.Range("F20").Value = 379133
$ws.Range("G20").Value = 416186
$ws.Range("H20").Value = 1661703

# Row 21: مالیات
$ws.Range("D21").Value = -36702
$ws.Range("E21").Value = -34279
$ws.Range("F21").Value = -75541
$ws.Range("G21").Value = -93669
$ws.Range("H21").Value = -295785

# Row 22: سود (زیان) خالص عملیات در حال تداوم
$ws.Range("D22").Value = 102607
$ws.Range("E22").Value = 157141
$ws.Range("F22").Value = 303592
$ws.Range("G22").Value = 322517
$ws.Range("H22").Value = 1365918

# Row 23: سود (زیان) عملیات متوقف شده پس از اثر مالیاتی -> now 0 instead of "-"
$ws.Range("D23").Value = 0
$ws.Range("E23").Value = 0
$ws.Range("F23").Value = 0
$ws.Range("G23").Value = 0
$ws.Range("H23").Value = 0

# Row 24: سود (زیان) خالص
$ws.Range("D24").Value = 102607
$ws.Range("E24").Value = 157141
$ws.Range("F24").Value = 303592
$ws.Range("G24").Value = 322517
$ws.Range("H24").Value = 1365918

# Row 25: سود هر سهم پس از کسر مالیات
$ws.Range("D25").Value = 513
$ws.Range("E25").Value = 100
$ws.Range("F25").Value = 193
$ws.Range("G25").Value = 188
$ws.Range("H25").Value = 795

# Row 26: سرمایه
$ws.Range("D26").Value = 200000
$ws.Range("E26").Value = 1574045
$ws.Range("F26").Value = 1574045
$ws.Range("G26").Value = 1718000
$ws.Range("H26").Value = 1718000

# Row 27: سود هر سهم بر اساس آخرین سرمایه
$ws.Range("D27").Value = 60
$ws.Range("E27").Value = 91
$ws.Range("F27").Value = 177
$ws.Range("G27").Value = 188
$ws.Range("H27").Value = 795
